# Update the subtitle on the title slide to reflect the final submission.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(5)
$shape.TextFrame.TextRange.Text = "Final Project by Nils Berzins"
